$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Multiply the Power values in F6:F12 by 1000 (unit conversion),
# matching the exact target values from the diff.
$ws.Range("F6").Value = 0.0267
$ws.Range("F7").Value = 0.28464
$ws.Range("F8").Value = 0.7695599999999999
$ws.Range("F9").Value = 10.88682
$ws.Range("F10").Value = 6.765119999999999
$ws.Range("F11").Value = 6.68448
$ws.Range("F12").Value = 11.6844
